$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.664.10"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "1.611.01"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.521"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.72"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.257"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0608"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.50%  "
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").Value = "1.842.61"
$ws.Range("D13").Value = "1.615.45"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.563"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.06%  "
$ws.Range("D16").Value = "29.710.88"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.62"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +13.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("D20").Value = "0.0₃0703"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.995"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.64%  "
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0481"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.01%  "
$ws.Range("E31").Value = "  +2.43%  "
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").Value = "1.454.09"
$ws.Range("E33").Value = "  +2.49%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.73%  "
$ws.Range("E36").Value = "  +1.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.63%  "
$ws.Range("E38").Value = "  -0.65%  "
$ws.Range("E39").Value = "  +2.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.554"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0506"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.822"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.65%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "69.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.73%  "
$ws.Range("B45").Value = "BitcoinSV"
$ws.Range("C45").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "53.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("E47").Value = "  +19.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.95%  "
$ws.Range("D49").Value = "1.752.72"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "87.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.09%  "
$ws.Range("E51").Value = "  -0.60%  "
